$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-09-04 08:45:22"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-09-04 08:45:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-09-04 08:45:22"
